$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.050.77"
$ws.Range("E2").Value = "  +0.81%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.645.28"
$ws.Range("E3").Value = "  +1.05%  "
$ws.Range("E4").Value = "  +1.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.50"
$ws.Range("E5").Value = "  +1.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.508"
$ws.Range("E6").Value = "  +1.07%  "
$ws.Range("E7").Value = "  +1.02%  "
$ws.Range("E8").Value = "  +0.61%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0640"
$ws.Range("E9").Value = "  +1.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.66"
$ws.Range("E10").Value = "  +0.40%  "
$ws.Range("E11").Value = "  +0.98%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.30"
$ws.Range("E12").Value = "  +1.64%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.872.11"
$ws.Range("E13").Value = "  +0.97%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.662.06"
$ws.Range("E14").Value = "  +1.68%  "
$ws.Range("E15").Value = "  +0.39%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₃0766"
$ws.Range("E16").Value = "  +1.29%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.28"
$ws.Range("E17").Value = "  +0.92%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.105.65"
$ws.Range("E18").Value = "  +1.02%  "
$ws.Range("E19").Value = "  +1.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "193.11"
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("E21").Value = "  -0.36%  "
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("E24").Value = "  +1.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.36"
$ws.Range("E25").Value = "  +1.74%  "
$ws.Range("E26").Value = "  +1.26%  "
$ws.Range("E27").Value = "  +4.01%  "
$ws.Range("E28").Value = "  +0.87%  "
$ws.Range("E29").Value = "  +0.54%  "
$ws.Range("E30").Value = "  +1.39%  "
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("E32").Value = "  -0.59%  "
$ws.Range("E33").Value = "  +1.35%  "
$ws.Range("E34").Value = "  +2.44%  "
$ws.Range("E35").Value = "  -3.05%  "
$ws.Range("E36").Value = "  +0.79%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.133.91"
$ws.Range("E37").Value = "  +0.20%  "
$ws.Range("E38").Value = "  -1.40%  "
$ws.Range("E39").Value = "  +0.31%  "
$ws.Range("E40").Value = "  +0.71%  "
$ws.Range("E41").Value = "  +0.81%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "99.53"
$ws.Range("E42").Value = "  +0.56%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.799"
$ws.Range("E43").Value = "  -0.56%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.781.44"
$ws.Range("E44").Value = "  +0.96%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0₆0118"
$ws.Range("E45").Value = "  +6.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.65"
$ws.Range("E46").Value = "  +1.15%  "
$ws.Range("E47").Value = "  +1.12%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.75"
$ws.Range("E48").Value = "  +2.01%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.45"
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.417"
$ws.Range("E50").Value = "  +0.72%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0958"
$ws.Range("E51").Value = "  -0.24%  "
